$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.013193762863988
$ws.Range("D2").Value = 1.015839054349307
$ws.Range("E2").Value = 1.015006082329285
$ws.Range("F2").Value = 1.024440664508684
$ws.Range("I2").Value = 1.025606812330986
$ws.Range("J2").Value = 1.018432134943767
$ws.Range("K2").Value = 1.018693014191807
$ws.Range("L2").Value = 1.017862533226702
$ws.Range("M2").Value = 1.027269158509428
$ws.Range("N2").Value = 1.009893179867788
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.015326193525176
$ws.Range("D3").Value = 1.017882829165986
$ws.Range("E3").Value = 1.016857263799404
$ws.Range("F3").Value = 1.026463788833949
$ws.Range("I3").Value = 1.02563770556681
$ws.Range("J3").Value = 1.020192331381777
$ws.Range("K3").Value = 1.020538611398353
$ws.Range("L3").Value = 1.019515883541158
$ws.Range("M3").Value = 1.029096070049915
$ws.Range("N3").Value = 1.010511492568047
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.016700364374441
$ws.Range("D4").Value = 1.01920009781476
$ws.Range("E4").Value = 1.018050135534325
$ws.Range("F4").Value = 1.027765535553417
$ws.Range("I4").Value = 1.025654747514573
$ws.Range("J4").Value = 1.021325661524972
$ws.Range("K4").Value = 1.021727303580986
$ws.Range("L4").Value = 1.020580357523859
$ws.Range("M4").Value = 1.030270504295453
$ws.Range("N4").Value = 1.010908673224233
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.017276745396144
$ws.Range("D5").Value = 1.01975266561617
$ws.Range("E5").Value = 1.018550457328437
$ws.Range("F5").Value = 1.028311059398515
$ws.Range("I5").Value = 1.02566120715142
$ws.Range("J5").Value = 1.021800792128579
$ws.Range("K5").Value = 1.022225733782121
$ws.Range("L5").Value = 1.02102660628689
$ws.Range("M5").Value = 1.030762420045761
$ws.Range("N5").Value = 1.011074961240664
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.01737344583738
$ws.Range("D6").Value = 1.019845373997367
$ws.Range("E6").Value = 1.018634396243218
$ws.Range("F6").Value = 1.028402554617162
$ws.Range("I6").Value = 1.025662250450633
$ws.Range("J6").Value = 1.021880491956386
$ws.Range("K6").Value = 1.022309347232147
$ws.Range("L6").Value = 1.021101460501895
$ws.Range("M6").Value = 1.030844909126277
$ws.Range("N6").Value = 1.011102841739707
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.016708071149521
$ws.Range("D7").Value = 1.019207485973731
$ws.Range("E7").Value = 1.018056825385526
$ws.Range("F7").Value = 1.02777283163009
$ws.Range("I7").Value = 1.025654836596532
$ws.Range("J7").Value = 1.0213320153985
$ws.Range("K7").Value = 1.021733968683787
$ws.Range("L7").Value = 1.020586325222746
$ws.Range("M7").Value = 1.030277084398754
$ws.Range("N7").Value = 1.010910897858182
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.013915618730264
$ws.Range("D8").Value = 1.016530850876518
$ws.Range("E8").Value = 1.015632742780939
$ws.Range("F8").Value = 1.025125928380495
$ws.Range("I8").Value = 1.025617863854005
$ws.Range("J8").Value = 1.019028184969715
$ws.Range("K8").Value = 1.019317905268105
$ws.Range("L8").Value = 1.018422414922443
$ws.Range("M8").Value = 1.027888182414642
$ws.Range("N8").Value = 1.010102749746004
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.008950146588167
$ws.Range("D9").Value = 1.011773137018571
$ws.Range("E9").Value = 1.011321915271883
$ws.Range("F9").Value = 1.020404131336477
$ws.Range("I9").Value = 1.025530088701005
$ws.Range("J9").Value = 1.014924136612974
$ws.Range("K9").Value = 1.01501682563427
$ws.Range("L9").Value = 1.014567155098706
$ws.Range("M9").Value = 1.023618438524257
$ws.Range("N9").Value = 1.008655969886184
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.005607572478285
$ws.Range("D10").Value = 1.008571706134578
$ws.Range("E10").Value = 1.008419857104194
$ws.Range("F10").Value = 1.01721565766169
$ws.Range("I10").Value = 1.025456291900837
$ws.Range("J10").Value = 1.012156483934058
$ws.Range("K10").Value = 1.012118252963713
$ws.Range("L10").Value = 1.011966979735253
$ws.Range("M10").Value = 1.020729740931483
$ws.Range("N10").Value = 1.007675552062586
$ws.Range("B11").Value = 1.019999999999999
$ws.Range("C11").Value = 1.004152070080325
$ws.Range("D11").Value = 1.007177983666579
$ws.Range("E11").Value = 1.007156154887436
$ws.Range("F11").Value = 1.015824958637285
$ws.Range("I11").Value = 1.02542069627375
$ws.Range("J11").Value = 1.010950162025065
$ws.Range("K11").Value = 1.010855330390328
$ws.Range("L11").Value = 1.010833587762245
$ws.Range("M11").Value = 1.019468501979665
$ws.Range("N11").Value = 1.007247108944062
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.003610166227729
$ws.Range("D12").Value = 1.006659130654608
$ws.Range("E12").Value = 1.006685659660439
$ws.Range("F12").Value = 1.015306841499203
$ws.Range("I12").Value = 1.025406925875131
$ws.Range("J12").Value = 1.010500857372336
$ws.Range("K12").Value = 1.010385013838931
$ws.Range("L12").Value = 1.010411436612878
$ws.Range("M12").Value = 1.018998422965191
$ws.Range("N12").Value = 1.007087365537009
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.003726464384762
$ws.Range("D13").Value = 1.006770479612945
$ws.Range("E13").Value = 1.006786632819329
$ws.Range("F13").Value = 1.015418050105905
$ws.Range("I13").Value = 1.025409904523278
$ws.Range("J13").Value = 1.010597290650125
$ws.Range("K13").Value = 1.010485953711809
$ws.Range("L13").Value = 1.010502042434882
$ws.Range("M13").Value = 1.019099329409692
$ws.Range("N13").Value = 1.007121658415814
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.00410730219742
$ws.Range("D14").Value = 1.007135119052572
$ws.Range("E14").Value = 1.007117286245754
$ws.Range("F14").Value = 1.015782162740211
$ws.Range("I14").Value = 1.025419569208786
$ws.Range("J14").Value = 1.010913047513225
$ws.Range("K14").Value = 1.01081647878838
$ws.Range("L14").Value = 1.010798716446792
$ws.Range("M14").Value = 1.019429677914884
$ws.Range("N14").Value = 1.007233916827147
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.004341779998811
$ws.Range("D15").Value = 1.007359630252192
$ws.Range("E15").Value = 1.007320865937265
$ws.Range("F15").Value = 1.01600629809943
$ws.Range("I15").Value = 1.02542545119822
$ws.Range("J15").Value = 1.011107432634816
$ws.Range("K15").Value = 1.011019964662021
$ws.Range("L15").Value = 1.010981352554649
$ws.Range("M15").Value = 1.019633003653952
$ws.Range("N15").Value = 1.007303002988979
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.005703994219459
$ws.Range("D16").Value = 1.008664042012205
$ws.Range("E16").Value = 1.008503572511613
$ws.Range("F16").Value = 1.017307738505443
$ws.Range("I16").Value = 1.025458577419319
$ws.Range("J16").Value = 1.012236373952719
$ws.Range("K16").Value = 1.012201901146182
$ws.Range("L16").Value = 1.012042038492712
$ws.Range("M16").Value = 1.020813222899502
$ws.Range("N16").Value = 1.007703902879983
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.006556265196287
$ws.Range("D17").Value = 1.009480235087914
$ws.Range("E17").Value = 1.009243530254226
$ws.Range("F17").Value = 1.018121376158442
$ws.Range("I17").Value = 1.025478380708179
$ws.Range("J17").Value = 1.012942387052339
$ws.Range("K17").Value = 1.012941179584229
$ws.Range("L17").Value = 1.012705348488431
$ws.Range("M17").Value = 1.02155073127031
$ws.Range("N17").Value = 1.00795431985342
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.00705259734141
$ws.Range("D18").Value = 1.009955587444754
$ws.Range("E18").Value = 1.009674453183094
$ws.Range("F18").Value = 1.018594988260295
$ws.Range("I18").Value = 1.02548958042486
$ws.Range("J18").Value = 1.013353431712695
$ws.Range("K18").Value = 1.013371635979825
$ws.Range("L18").Value = 1.013091524663187
$ws.Range("M18").Value = 1.021979904981986
$ws.Range("N18").Value = 1.008100006865976
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.007221702199079
$ws.Range("D19").Value = 1.010117549467267
$ws.Range("E19").Value = 1.00982127213742
$ws.Range("F19").Value = 1.018756314497618
$ws.Range("I19").Value = 1.025493339714151
$ws.Range("J19").Value = 1.013493459306006
$ws.Range("K19").Value = 1.013518283972265
$ws.Range("L19").Value = 1.013223079378375
$ws.Range("M19").Value = 1.022126073118413
$ws.Range("N19").Value = 1.008149618799085
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.006464905884827
$ws.Range("D20").Value = 1.009392739966967
$ws.Range("E20").Value = 1.009164210571543
$ws.Range("F20").Value = 1.01803418102897
$ws.Range("I20").Value = 1.025476292340809
$ws.Range("J20").Value = 1.012866717473723
$ws.Range("K20").Value = 1.012861940065976
$ws.Range("L20").Value = 1.012634256445203
$ws.Range("M20").Value = 1.02147170753457
$ws.Range("N20").Value = 1.007927491562098
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.003995190224707
$ws.Range("D21").Value = 1.00702777422648
$ws.Range("E21").Value = 1.00701994769722
$ws.Range("F21").Value = 1.01567498372072
$ws.Range("I21").Value = 1.025416738356609
$ws.Range("J21").Value = 1.010820099028342
$ws.Range("K21").Value = 1.01071918110818
$ws.Range("L21").Value = 1.010711385598358
$ws.Range("M21").Value = 1.019332442923583
$ws.Range("N21").Value = 1.007200876200891
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.00243503924645
$ws.Range("D22").Value = 1.00553408140702
$ws.Range("E22").Value = 1.005665382969716
$ws.Range("F22").Value = 1.014182678456709
$ws.Range("I22").Value = 1.025376119566924
$ws.Range("J22").Value = 1.009526215375146
$ws.Range("K22").Value = 1.009364919546081
$ws.Range("L22").Value = 1.009495678880874
$ws.Range("M22").Value = 1.017978135306703
$ws.Range("N22").Value = 1.00674054340692
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.003262814907603
$ws.Range("D23").Value = 1.006326568576381
$ws.Range("E23").Value = 1.006384079935383
$ws.Range("F23").Value = 1.014974641810342
$ws.Range("I23").Value = 1.025397953830631
$ws.Range("J23").Value = 1.010212811668405
$ws.Range("K23").Value = 1.010083517133198
$ws.Range("L23").Value = 1.010140796067329
$ws.Range("M23").Value = 1.018696969382584
$ws.Range("N23").Value = 1.006984908632195
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.006506189643186
$ws.Range("D24").Value = 1.009432277466876
$ws.Range("E24").Value = 1.00920005382757
$ws.Range("F24").Value = 1.018073583747463
$ws.Range("I24").Value = 1.025477237069366
$ws.Range("J24").Value = 1.012900911648704
$ws.Range("K24").Value = 1.012897747316312
$ws.Range("L24").Value = 1.01266638210699
$ws.Range("M24").Value = 1.021507418055492
$ws.Range("N24").Value = 1.007939615276504
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.01023936779762
$ws.Range("D25").Value = 1.013008195302354
$ws.Range("E25").Value = 1.012441202297609
$ws.Range("F25").Value = 1.021631847346842
$ws.Range("I25").Value = 1.025555468655154
$ws.Range("J25").Value = 1.015990573526646
$ws.Range("K25").Value = 1.016134118743947
$ws.Range("L25").Value = 1.015568997755284
$ws.Range("M25").Value = 1.024729579113322
$ws.Range("N25").Value = 1.009032752944889

Write-Host "done"
